$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 2.35
$ws.Range("I7").Value = 3.7
$ws.Range("J7").Value = 3.4
$ws.Range("AA7").Value = 5
$ws.Range("AD7").Value = 23
$ws.Range("AL7").Value = 6.5
$ws.Range("AM7").Value = 15
$ws.Range("AN7").Value = 15
# Row 8
$ws.Range("G8").Value = 1.62
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 6.25
$ws.Range("J8").Value = 2.3
$ws.Range("K8").Value = 1.95
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.88
$ws.Range("S8").Value = 2.6
$ws.Range("T8").Value = 1.48
$ws.Range("U8").Value = 5.5
$ws.Range("V8").Value = 1.14
$ws.Range("W8").Value = 1.57
$ws.Range("X8").Value = 2.25
$ws.Range("Y8").Value = 2.5
$ws.Range("Z8").Value = 1.5
$ws.Range("AA8").Value = 4.75
$ws.Range("AC8").Value = 9.5
$ws.Range("AD8").Value = 11
$ws.Range("AE8").Value = 17
$ws.Range("AG8").Value = 6.5
$ws.Range("AI8").Value = 26
$ws.Range("AL8").Value = 11
$ws.Range("AR8").Value = 4
$ws.Range("AS8").Value = 1.25
# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 2.8
$ws.Range("I10").Value = 5.75
$ws.Range("AL10").Value = 9.5
$ws.Range("AM10").Value = 26
# Row 11
$ws.Range("G11").Value = 1.95
$ws.Range("I11").Value = 4.33
$ws.Range("J11").Value = 2.75
$ws.Range("L11").Value = 6
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("U11").Value = 6.5
$ws.Range("V11").Value = 1.11
$ws.Range("AA11").Value = 4.75
$ws.Range("AB11").Value = 7.5
$ws.Range("AI11").Value = 26
$ws.Range("AJ11").Value = 126
$ws.Range("AL11").Value = 8
$ws.Range("AM11").Value = 21
$ws.Range("AP11").Value = 51
$ws.Range("AQ11").Value = 67
# Row 16
$ws.Range("G16").Value = 2.55
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 2.75
$ws.Range("J16").Value = 3.2
$ws.Range("L16").Value = 3.4
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 9.5
$ws.Range("O16").Value = 1.3
$ws.Range("P16").Value = 3.4
$ws.Range("Q16").Value = 1.5
$ws.Range("R16").Value = 2.55
$ws.Range("S16").Value = 2.03
$ws.Range("T16").Value = 1.83
$ws.Range("AA16").Value = 8.5
$ws.Range("AB16").Value = 12
$ws.Range("AD16").Value = 23
$ws.Range("AG16").Value = 9.5
$ws.Range("AH16").Value = 6
$ws.Range("AI16").Value = 13
$ws.Range("AJ16").Value = 41
$ws.Range("AL16").Value = 9
$ws.Range("AM16").Value = 13
$ws.Range("AN16").Value = 11
$ws.Range("AO16").Value = 29
$ws.Range("AP16").Value = 23
$ws.Range("AR16").Value = 2.75
$ws.Range("AS16").Value = 1.44
# Row 18
$ws.Range("G18").Value = 2.3
$ws.Range("H18").Value = 3
$ws.Range("J18").Value = 3.2
$ws.Range("K18").Value = 1.83
$ws.Range("M18").Value = 1.14
$ws.Range("N18").Value = 5.5
$ws.Range("O18").Value = 1.62
$ws.Range("P18").Value = 2.2
$ws.Range("Q18").Value = 2.1
$ws.Range("R18").Value = 1.78
$ws.Range("S18").Value = 2.88
$ws.Range("T18").Value = 1.4
$ws.Range("U18").Value = 6
$ws.Range("V18").Value = 1.13
$ws.Range("W18").Value = 1.67
$ws.Range("X18").Value = 2.1
$ws.Range("Y18").Value = 2.38
$ws.Range("Z18").Value = 1.53
$ws.Range("AC18").Value = 11
$ws.Range("AG18").Value = 5.5
$ws.Range("AL18").Value = 7
$ws.Range("AP18").Value = 41
$ws.Range("AR18").Value = 4.6
$ws.Range("AS18").Value = 1.18
# Row 20
$ws.Range("G20").Value = 1.67
$ws.Range("H20").Value = 3.75
$ws.Range("I20").Value = 4.75
$ws.Range("J20").Value = 2.3
$ws.Range("K20").Value = 2.25
$ws.Range("L20").Value = 5
$ws.Range("AB20").Value = 8.5
$ws.Range("AD20").Value = 13
$ws.Range("AI20").Value = 15
$ws.Range("AK20").Value = 201
$ws.Range("AL20").Value = 13
$ws.Range("AM20").Value = 23
$ws.Range("AN20").Value = 15
# Row 46
$ws.Range("G46").Value = 2.52
$ws.Range("I46").Value = 2.72
$ws.Range("J46").Value = 3.15
$ws.Range("K46").Value = 1.98
$ws.Range("L46").Value = 3.4
$ws.Range("O46").Value = 1.39
$ws.Range("P46").Value = 2.55
$ws.Range("S46").Value = 2.22
$ws.Range("T46").Value = 1.52
$ws.Range("U46").Value = 3.95
$ws.Range("V46").Value = 1.17
$ws.Range("W46").Value = 1.42
$ws.Range("X46").Value = 2.45
$ws.Range("Y46").Value = 1.87
$ws.Range("Z46").Value = 1.75
$ws.Range("AA46").Value = 7
$ws.Range("AB46").Value = 12
$ws.Range("AC46").Value = 9.75
$ws.Range("AD46").Value = 29
$ws.Range("AE46").Value = 24
$ws.Range("AF46").Value = 37
$ws.Range("AG46").Value = 7.2
$ws.Range("AH46").Value = 5.8
$ws.Range("AI46").Value = 15.5
$ws.Range("AK46").Value = 800
$ws.Range("AL46").Value = 7.3
$ws.Range("AM46").Value = 13
$ws.Range("AN46").Value = 10.5
$ws.Range("AO46").Value = 35
$ws.Range("AP46").Value = 27
